# Fruta / hortaliza, semanal
# A new daily price observation is inserted at row 70 (pushing every
# subsequent record down by one row); the previously-last record (old
# row 155) ends up duplicated down into the newly created row 156.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 70 - this shifts rows 70:155
# down to 71:156 and copies the formatting (incl. the date number
# format on column D) from the row above, matching native Excel
# "Insert" behaviour.
$ws.Rows("70:70").Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A70").Value = 7
$ws.Range("B70").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C70").Value = "Ñuble"
$ws.Range("D70").Value = 44483
$ws.Range("E70").Value = 16
$ws.Range("F70").Value = 100112043
$ws.Range("G70").Value = "Pepino ensalada"
$ws.Range("H70").Value = "Sin especificar"
$ws.Range("I70").Value = "Primera"
$ws.Range("J70").Value = 160
$ws.Range("K70").Value = 16000
$ws.Range("L70").Value = 17000
$ws.Range("M70").Value = 16500
$ws.Range("N70").Value = "$/caja 80 unidades"
$ws.Range("O70").Value = "Región del Maule"
$ws.Range("P70").Value = 206
$ws.Range("Q70").Value = 80
$ws.Range("R70").Value = "Hortaliza"
